$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 132, shifting existing data down
$ws.Rows.Item(132).Insert()
$ws.Rows.Item(132).Insert()

# Populate new row 132 (Primera quality, date 2022-04-29)
$ws.Range("A132").Value = 11
$ws.Range("B132").Value = "Vega Monumental Concepción"
$ws.Range("C132").Value = "Bíobío"
$ws.Range("D132").Value = 44680
$ws.Range("E132").Value = 8
$ws.Range("F132").Value = 100114013
$ws.Range("G132").Value = "Zanahoria"
$ws.Range("H132").Value = "Sin especificar"
$ws.Range("I132").Value = "Primera"
$ws.Range("J132").Value = 1000
$ws.Range("K132").Value = 6000
$ws.Range("L132").Value = 6500
$ws.Range("M132").Value = 6250
$ws.Range("N132").Value = "$/saco 20 kilos"
$ws.Range("O132").Value = "Región de Ñuble"
$ws.Range("P132").Value = 312
$ws.Range("Q132").Value = 20
$ws.Range("R132").Value = "Hortaliza"

# Populate new row 133 (Segunda quality, date 2022-04-29)
$ws.Range("A133").Value = 11
$ws.Range("B133").Value = "Vega Monumental Concepción"
$ws.Range("C133").Value = "Bíobío"
$ws.Range("D133").Value = 44680
$ws.Range("E133").Value = 8
$ws.Range("F133").Value = 100114013
$ws.Range("G133").Value = "Zanahoria"
$ws.Range("H133").Value = "Sin especificar"
$ws.Range("I133").Value = "Segunda"
$ws.Range("J133").Value = 500
$ws.Range("K133").Value = 5000
$ws.Range("L133").Value = 5000
$ws.Range("M133").Value = 5000
$ws.Range("N133").Value = "$/saco 20 kilos"
$ws.Range("O133").Value = "Región de Ñuble"
$ws.Range("P133").Value = 250
$ws.Range("Q133").Value = 20
$ws.Range("R133").Value = "Hortaliza"
